# The deck originally ships an "Integral" theme (ppt/theme/theme1.xml,
# used by the slide master / all slides) and an unused "Office Theme"
# (ppt/theme/theme2.xml, only linked from the notes master). The commit
# swaps the two themes' content, so the slide master ends up using the
# default "Office" palette instead of "Integral".
#
# Font scheme and format scheme (fills/lines/effects) are byte-identical
# between the two themes, so only the 12 theme colours actually change.
# dk1/lt1 (black/white) are identical between the two palettes too, so
# only dk2, lt2 and the six accents + the two hyperlink colours move.

function Get-RGBLong([string]$hex) {
    $r = [Convert]::ToInt32($hex.Substring(0, 2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2, 2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4, 2), 16)
    return $r + ($g * 256) + ($b * 65536)
}

$p = $ppt.ActivePresentation
$theme = $p.SlideMaster.Theme
$scheme = $theme.ThemeColorScheme

# msoThemeColorSchemeIndex order:
# 1 Dark1, 2 Light1, 3 Dark2, 4 Light2,
# 5 Accent1, 6 Accent2, 7 Accent3, 8 Accent4, 9 Accent5, 10 Accent6,
# 11 Hyperlink, 12 FollowedHyperlink
$scheme.Item(3).RGB  = Get-RGBLong "44546A"   # dk2
$scheme.Item(4).RGB  = Get-RGBLong "E7E6E6"   # lt2
$scheme.Item(5).RGB  = Get-RGBLong "5B9BD5"   # accent1
$scheme.Item(6).RGB  = Get-RGBLong "ED7D31"   # accent2
$scheme.Item(7).RGB  = Get-RGBLong "A5A5A5"   # accent3
$scheme.Item(8).RGB  = Get-RGBLong "FFC000"   # accent4
$scheme.Item(9).RGB  = Get-RGBLong "4472C4"   # accent5
$scheme.Item(10).RGB = Get-RGBLong "70AD47"   # accent6
$scheme.Item(11).RGB = Get-RGBLong "0563C1"   # hlink
$scheme.Item(12).RGB = Get-RGBLong "954F72"   # folHlink
